# Allow full switch to bio options in IND by 2050
# Set the "by 2050" mitigation-option switch (column T) from 0.5 (50%) to 1
# (100%) for every mitigation row on the IND_Shares sheet, so that the
# associated O-column impact formulas (O = <base>*T) recalc to double their
# previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IND_Shares")

$rows = @(12,13,14,24,26,27,34,35,36,44,45,46,53,54,55,64,65,66,73,74,75,87,88,89,98,99,100,107,108,109,116,117,118,125,126,127,135,136,137,144,145,146)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 20).Value = 1
}

$ws.Activate()
$ws.Range("W81").Select()
